$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web")

# --- Add new "Progress Bar" widget section (rows 149-152) ---

# Row 149: category header
$ws.Range("A149").Value = "progressBar"

# Row 150: nav link
$ws.Range("A150").Value = "progressBarNav"
$ws.Range("B150").Value = "//div[@class='element-list collapse show']//li[@id='item-4']"
$ws.Range("C150").Value = "By.xpath"

# Row 151: start/stop button
$ws.Range("A151").Value = "startProgressButton"
$ws.Range("B151").Value = "//button[@id='startStopButton']"
$ws.Range("C151").Value = "By.xpath"

# Row 152: progress bar value
$ws.Range("A152").Value = "progressBarValue"
$ws.Range("B152").Value = '//*[@id="progressBar"]/div'
$ws.Range("C152").Value = "By.xpath"

# --- Update view state (scroll position / active selection) ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 145
$ws.Range("A155").Select()
